$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 ("Summary of requirements" list) - Content Placeholder (shape 2)
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
$tr5.Text = "Summary of requirements"
[void]$tr5.InsertAfter("`rCalhoun- Billing")
[void]$tr5.InsertAfter("`rDante-")
[void]$tr5.InsertAfter("`rShaun-")
[void]$tr5.InsertAfter("`rKincaid-")
[void]$tr5.InsertAfter("`r")

# ---------------------------------------------------------------------------
# 2) Add slide 7 as a duplicate of slide 6 *before* slide 6 gets its own new
#    text, so the new slide starts from the same "empty placeholder" state
#    that slide 6 currently has (avoids picking up layout-inherited bodyPr
#    overrides such as anchor="ctr").
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$dup = $s6.Duplicate()
$s7 = $p.Slides.Item($p.Slides.Count)

# ---------------------------------------------------------------------------
# 3) Slide 6 ("Cost, almost last slide") - Content Placeholder (shape 2)
# ---------------------------------------------------------------------------
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Cost, almost last slide"

# ---------------------------------------------------------------------------
# 4) New slide 7 ("conclusion") - Title (shape 1); content stays blank
# ---------------------------------------------------------------------------
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "conclusion"
